# Update of league bases (Colombia Primera B) - 04-04-2024 23:22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Three pairs of rows had their match-data columns swapped (id/Div/
# Div Original Name/Date stay put in A/C/D/E; everything else - the two
# team columns, score, result, and all odds columns B and F..AC - trade
# places between the two rows).
# ---------------------------------------------------------------------------
$swapCols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
$swapPairs = @(
    @(176,177),
    @(188,189),
    @(227,228)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $swapCols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

# ---------------------------------------------------------------------------
# Rows 243-245 get refreshed match data (new/updated fixtures), and the
# two trailing placeholder rows (246-247) are removed entirely - the sheet
# shrinks from 247 to 245 data rows.
# ---------------------------------------------------------------------------
$ws.Range("B243").Value2 = 7657923
$ws.Range("E243").Value2 = 45387.72916666666
$ws.Range("F243").Value2 = "Tigres FC"
$ws.Range("G243").Value2 = "Atletico Huila"
$ws.Range("K243").Value2 = 2.75
$ws.Range("L243").Value2 = 3.4
$ws.Range("M243").Value2 = 2.35
$ws.Range("N243").Value2 = 2.875
$ws.Range("O243").Value2 = 3.4
$ws.Range("P243").Value2 = 2.25
$ws.Range("Q243").Value2 = 0.25
$ws.Range("R243").Value2 = 1.8
$ws.Range("S243").Value2 = 2
$ws.Range("T243").Value2 = 2.25
$ws.Range("U243").Value2 = 2
$ws.Range("V243").Value2 = 1.8

$ws.Range("B244").Value2 = 7658165
$ws.Range("E244").Value2 = 45388.72916666666
$ws.Range("F244").Value2 = "Barranquilla FC"
$ws.Range("G244").Value2 = "Deportes Quindio"
$ws.Range("K244").Value2 = 2.6
$ws.Range("L244").Value2 = 3
$ws.Range("M244").Value2 = 2.75
$ws.Range("N244").Value2 = 3
$ws.Range("O244").Value2 = 3.1
$ws.Range("P244").Value2 = 2.55
$ws.Range("Q244").Value2 = 0
$ws.Range("R244").Value2 = 2.05
$ws.Range("S244").Value2 = 1.75
$ws.Range("T244").Value2 = 2.25
$ws.Range("U244").Value2 = 1.9
$ws.Range("V244").Value2 = 1.9

$ws.Range("B245").Value2 = 7658162
$ws.Range("E245").Value2 = 45388.72916666666
$ws.Range("F245").Value2 = "Real San Andres"
$ws.Range("G245").Value2 = "Bogota FC"
$ws.Range("K245").Value2 = 1.8
$ws.Range("L245").Value2 = 3.3
$ws.Range("M245").Value2 = 4.333
$ws.Range("N245").Value2 = 1.909
$ws.Range("O245").Value2 = 3.4
$ws.Range("P245").Value2 = 4.2
$ws.Range("Q245").Value2 = -0.5
$ws.Range("R245").Value2 = 1.85
$ws.Range("S245").Value2 = 1.95
$ws.Range("T245").Value2 = 2.5
$ws.Range("U245").Value2 = 1.975
$ws.Range("V245").Value2 = 1.825

# Drop the two now-obsolete rows (old 246 and 247); deleting 246 twice
# pulls what was 247 up into its place before removing it too.
$ws.Rows.Item(246).Delete()
$ws.Rows.Item(246).Delete()
